$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update Version / Date / Contact values ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(3, 2).Value = "2.0.0"
$meta.Cells.Item(8, 2).Value = "2024-06-03T10:45:43+02:00"
$meta.Cells.Item(10, 2).Value = "Kommunernes Landsforening (http://kl.dk)"

# --- Sheet "Include from FSIII": insert 5 new concept rows, shifting the
#     existing J1..J4 rows down, clear the old System-URI-looking values
#     out of column B for the (now bare) concept rows, and append the
#     System URI row at the bottom ---
$concept = $wb.Worksheets.Item("Include from FSIII")

# Extend the existing row style (s="2") down through the new rows 9-13 by
# copying formatting from row 8 (last pre-existing data row).
$concept.Range("A8:B8").Copy()
$concept.Range("A9:B13").PasteSpecial(-4122)

# Rows 2-6 become the five new UUID concept codes (column B stays blank).
$concept.Cells.Item(2, 1).Value = "43c2b7f0-5e55-4627-8fcf-bdaf5a9d84ac"
$concept.Cells.Item(3, 1).Value = "86b53158-6d05-412e-ad55-2e1fa26359b3"
$concept.Cells.Item(4, 1).Value = "1c850a09-aa49-4fae-9354-f932f13e030b"
$concept.Cells.Item(5, 1).Value = "462f9352-0129-4d8e-8c75-a6dfed78ddcf"
$concept.Cells.Item(6, 1).Value = "4571f168-a92a-4caf-8dc8-35f45c2a1cb4"

# Rows 7-11 hold the concept codes that used to occupy rows 2-6
# (J1, J5, J2, J3, J4); their column B values are cleared.
$concept.Cells.Item(7, 1).Value = "J1"
$concept.Cells.Item(7, 2).ClearContents()

$concept.Cells.Item(8, 1).Value = "J5"
$concept.Cells.Item(8, 2).ClearContents()

$concept.Cells.Item(9, 1).Value = "J2"
$concept.Cells.Item(10, 1).Value = "J3"
$concept.Cells.Item(11, 1).Value = "J4"

# Row 12 stays blank in both columns.
$concept.Cells.Item(12, 1).ClearContents()
$concept.Cells.Item(12, 2).ClearContents()

# Row 13 is the new System URI row.
$concept.Cells.Item(13, 1).Value = "System URI"
$concept.Cells.Item(13, 2).Value = "urn:oid:1.2.208.176.2.21"
